# "Updated_Format_1 (Text Wrapper+Coloana B mai mica)"
#
# 1) Name label gets two leading spaces.
# 2) The small gray "count" font (column A number cells) changes color
#    FF808080 -> FFE3E3E3 and becomes centered (horizontal+vertical).
# 3) Colored "word" header cells (column A) become centered + wrap text;
#    the sentence/filler cells in column B get wrap text turned on.
# 4) Column B gets narrower (250.71 -> 200.71 characters).
# 5) Each block gets a new "remaining sentences" counter in column A,
#    placed on the first data row below the block header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Name label: "Numele și Prenumele:" -> "  Numele și Prenumele:" ---
$ws.Range("A3").Value = "  Numele și Prenumele:"

# --- 4) Narrower column B ---
$ws.Columns(2).ColumnWidth = 200.7109375

# --- block layout: (wordRow, numberRow, firstDataRow, lastDataRow) ---
$blocks = @(
    @(5, 6, 7, 15),
    @(17, 18, 19, 27),
    @(29, 30, 31, 39),
    @(41, 42, 43, 51),
    @(53, 54, 55, 63),
    @(65, 66, 67, 75),
    @(77, 78, 79, 87),
    @(89, 90, 91, 99),
    @(101, 102, 103, 111),
    @(113, 114, 115, 123),
    @(125, 126, 127, 135),
    @(137, 138, 139, 147)
)

# --- 5) counters placed at the first data row of each block ---
$counters = @{
    7   = 10
    19  = 10
    31  = 10
    43  = 10
    55  = 3
    67  = 10
    79  = 10
    91  = 6
    103 = 8
    115 = 4
    127 = 10
    139 = 10
}

# column-A number cells whose formatting changes (the header A1 plus every
# "numberRow" A-cell, plus the brand-new counter cells)
$numberCells = @("A1")

foreach ($block in $blocks) {
    $wordRow = $block[0]
    $numberRow = $block[1]
    $firstData = $block[2]
    $lastData = $block[3]

    $numberCells += "A$numberRow"

    # --- 3a) word header cell (column A) gets centered + wrap text ---
    $wordCell = $ws.Range("A$wordRow")
    $wordCell.HorizontalAlignment = -4108   # xlCenter
    $wordCell.WrapText = $true

    # --- 3b) every sentence/filler cell in column B of the block wraps ---
    $dataRange = $ws.Range("B$firstData`:B$lastData")
    $dataRange.WrapText = $true

    # --- 5) drop the counter value onto the first data row (col A) ---
    if ($counters.ContainsKey($firstData)) {
        $ws.Range("A$firstData").Value = $counters[$firstData]
        $numberCells += "A$firstData"
    }
}

# --- 2) recolor + center every column-A "number" cell ---
foreach ($addr in $numberCells) {
    $cell = $ws.Range($addr)
    $cell.Font.Color = 14935011   # RGB(0xE3,0xE3,0xE3) == FFE3E3E3
    $cell.HorizontalAlignment = -4108  # xlCenter
    $cell.VerticalAlignment = -4108    # xlCenter
}

Write-Host "edit applied"
